$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI edge-weight values per Dr Hou's advice (ligand/receptor expressing
# cell counts changed from 1 to 3, with all derived average/total/specificity
# columns recomputed accordingly) for data rows 2-10.
$cols = @("E", "G", "H", "I", "J", "K", "M", "N", "O", "P", "Q", "R", "S", "T")
$rowData = @{
    2 = @(3, 60.90318633333334, 182.709559, 0.4799022665420342, 0.4799022665420342, 3, 73.91316300000001, 221.739489, 0.6096331558809399, 0.6096331558809398, 4501.54713867504, 40513.92424807536, 0.2925643332664363, 0.2925643332664362)
    3 = @(3, 60.90318633333334, 182.709559, 0.4799022665420342, 0.4799022665420342, 3, 31.40056566666667, 94.201697, 0.2589907557307037, 0.2589907557307037, 1912.394501769069, 17211.55051592162, 0.124290250688599, 0.124290250688599)
    4 = @(3, 60.90318633333334, 182.709559, 0.4799022665420342, 0.4799022665420342, 3, 15.928304, 47.784912, 0.1313760883883564, 0.1313760883883564, 970.0844664859786, 8730.760198373808, 0.06304768258699887, 0.06304768258699887)
    5 = @(3, 51.42568199999999, 154.277046, 0.405221842009972, 0.405221842009972, 3, 73.91316300000001, 221.739489, 0.6096331558809399, 0.6096331558809398, 3801.034816052166, 34209.31334446949, 0.2470366703764268, 0.2470366703764268)
    6 = @(3, 51.42568199999999, 154.277046, 0.405221842009972, 0.405221842009972, 3, 31.40056566666667, 94.201697, 0.2589907557307037, 0.2589907557307037, 1614.795504594118, 14533.15954134706, 0.1049487111007504, 0.1049487111007504)
    7 = @(3, 51.42568199999999, 154.277046, 0.405221842009972, 0.405221842009972, 3, 15.928304, 47.784912, 0.1313760883883564, 0.1313760883883564, 819.1238963033279, 7372.115066729951, 0.05323646053279468, 0.05323646053279468)
    8 = @(3, 14.57860966666667, 43.735829, 0.1148758914479938, 0.1148758914479938, 3, 73.91316300000001, 221.739489, 0.6096331558809399, 0.6096331558809398, 1077.551152605709, 9697.96037345138, 0.07003215223807675, 0.07003215223807674)
    9 = @(3, 14.57860966666667, 43.735829, 0.1148758914479938, 0.1148758914479938, 3, 31.40056566666667, 94.201697, 0.2589907557307037, 0.2589907557307037, 457.7765901668681, 4119.989311501812, 0.0297517939413542, 0.0297517939413542)
    10 = @(3, 14.57860966666667, 43.735829, 0.1148758914479938, 0.1148758914479938, 3, 15.928304, 47.784912, 0.1313760883883564, 0.1313760883883564, 232.2125266680053, 2089.912740012048, 0.01509194526856288, 0.01509194526856288)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
